$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.905.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5089"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06342"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.641.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7701"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.892.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "195.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.414"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.892"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.011"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "

# Row 23
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.865"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1194"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.810"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.233"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04897"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.239"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.165"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.524"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.367"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8888"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.582"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.139.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5395"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

# Row 41
$ws.Range("E41").Value = "  -1.03%  "

# Row 42
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₈127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.50%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8129"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.435"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.09%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.768.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4527"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9955"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05051"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "

Write-Output "Applied all crypto list updates"
